$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The dataset was rerun; 43 new observations (sheet rows 878..920,
# i.e. index 876..918 in column A) are appended below the previous last
# row (877).
$aValues = @(876,877,878,879,880,881,882,883,884,885,886,887,888,889,890,891,892,893,894,895,896,897,898,899,900,901,902,903,904,905,906,907,908,909,910,911,912,913,914,915,916,917,918)
$bValues = @(0.30726,0.30726,0.32464,0.31037,0.31595,0.29919,0.30478,0.30168,0.29919,0.2905,0.2905,0.2905,0.28181,0.2396,0.25388,0.23712,0.2396,0.23402,0.22222,0.21664,0.21167,0.20546,0.20608,0.17194,0.18063,0.20732,0.19429,0.19429,0.21788,0.19491,0.2067,0.20608,0.21167,0.21415,0.22222,0.21353,0.2185,0.20857,0.19181,0.21664,0.22843,0.26505,0.27933)

$lastExistingRow = 877
$startRow = $lastExistingRow + 1
$endRow = $startRow + $aValues.Length - 1

# Copy the per-column formatting of the last existing data row down across
# the whole new block (column A carries the bold/bordered "index" style,
# column B keeps the default style), then fill in the new values.
$ws.Range("A$lastExistingRow").Copy()
$ws.Range("A$($startRow):A$endRow").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("B$lastExistingRow").Copy()
$ws.Range("B$($startRow):B$endRow").PasteSpecial(-4122)  # xlPasteFormats

for ($i = 0; $i -lt $aValues.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $aValues[$i]
    $ws.Cells.Item($r, 2).Value = $bValues[$i]
}
